$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

for ($r = 84; $r -le 124; $r++) {
    $r2 = $r + 41

    # Snapshot both rows' original values before writing anything back,
    # since row r and row r2 both change as part of the same swap.
    $rowVals = @{}
    $row2Vals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range($col + $r).Value2
        $row2Vals[$col] = $ws.Range($col + $r2).Value2
    }
    $kText = $ws.Range("K" + $r).Text
    $k2Text = $ws.Range("K" + $r2).Text

    # Each row takes on the other's values, rescaled by x100 (fraction -> percent),
    # and the Private1/Public1 sector label swaps along with it.
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = $row2Vals[$col] * 100
        $ws.Range($col + $r2).Value = $rowVals[$col] * 100
    }
    $ws.Range("K" + $r).Value = $k2Text
    $ws.Range("K" + $r2).Value = $kText
}

# Match the author's saved selection/scroll state as closely as possible.
$ws.Range("C128").Select()
